# Update gh-pages to output generated at 456a3b4
# Applies the upstream data refresh to all four sheets:
#   1 展览     (Exhibition)
#   2 演出     (Performance)
#   3 本地生活 (Local life)
#   4 全部类型 (All types)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Simple "want to go" counter bumps (rows unaffected by the later insert)
$ws1.Cells.Item(2,6).Value  = 245
$ws1.Cells.Item(3,6).Value  = 583
$ws1.Cells.Item(6,6).Value  = 3183
$ws1.Cells.Item(7,6).Value  = 2759
$ws1.Cells.Item(8,6).Value  = 532
$ws1.Cells.Item(11,6).Value = 355
$ws1.Cells.Item(12,6).Value = 291
$ws1.Cells.Item(14,6).Value = 5735
$ws1.Cells.Item(20,6).Value = 465
$ws1.Cells.Item(21,6).Value = 1248

# A brand-new event ("上海·运动番ONLY") was added on 2024.04.13; it sorts
# before the existing row 23, so insert a fresh row there and push
# everything else down by one (old 23-26 -> new 24-27).
$ws1.Rows.Item(23).Insert()

# Match the bordered/bold/centered look of the other "序号" cells in col A.
$ws1.Cells.Item(22,1).Copy()
$ws1.Cells.Item(23,1).PasteSpecial(-4122)   # xlPasteFormats

$ws1.Cells.Item(23,1).Value = 22
# Force text format first so the pure "YYYY.MM.DD" string isn't
# auto-coerced into a date serial number by the COM layer (every other
# row in this column stores the start date as literal text).
$ws1.Cells.Item(23,2).NumberFormat = "@"
$ws1.Cells.Item(23,2).Value = "2024.04.13"
$ws1.Cells.Item(23,3).Value = "上海·运动番ONLY"
$ws1.Cells.Item(23,4).Value = "少年村路6号 YC篮羽联盟(大场店)"
$ws1.Cells.Item(23,5).Value = "2024.04.13 10:00-04.13 17:00"
$ws1.Cells.Item(23,6).Value = 1
$ws1.Cells.Item(23,7).Value = 60
$ws1.Cells.Item(23,8).Value = "https://show.bilibili.com/platform/detail.html?id=81901"
$ws1.Cells.Item(23,9).Value = "//i0.hdslb.com/bfs/openplatform/202402/2oiNlCAr1708325440584.jpeg"

# Sequence numbers and "want to go" counters for the rows that shifted down.
$ws1.Cells.Item(24,1).Value = 23
$ws1.Cells.Item(24,6).Value = 336

$ws1.Cells.Item(25,1).Value = 24
$ws1.Cells.Item(25,6).Value = 128

$ws1.Cells.Item(26,1).Value = 25

$ws1.Cells.Item(27,1).Value = 26
$ws1.Cells.Item(27,6).Value = 45

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performance)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(6,6).Value  = 243
$ws2.Cells.Item(8,6).Value  = 338
$ws2.Cells.Item(9,6).Value  = 62
$ws2.Cells.Item(12,6).Value = 15
$ws2.Cells.Item(13,6).Value = 635
$ws2.Cells.Item(19,6).Value = 623
$ws2.Cells.Item(24,6).Value = 292
$ws2.Cells.Item(25,6).Value = 4034
$ws2.Cells.Item(30,6).Value = 62

# Tickets for this show are no longer marked "不可售" (unavailable) - a
# price has now been published.
$ws2.Cells.Item(33,7).Value = 288

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(2,6).Value  = 1794
$ws3.Cells.Item(3,6).Value  = 85
$ws3.Cells.Item(5,6).Value  = 2568
$ws3.Cells.Item(6,6).Value  = 1119
$ws3.Cells.Item(10,6).Value = 403
$ws3.Cells.Item(13,6).Value = 423

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) - aggregated view, kept in sync by hand
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(2,6).Value  = 1794
$ws4.Cells.Item(3,6).Value  = 85
$ws4.Cells.Item(5,6).Value  = 2568
$ws4.Cells.Item(6,6).Value  = 1119
$ws4.Cells.Item(8,6).Value  = 403
$ws4.Cells.Item(11,6).Value = 245
$ws4.Cells.Item(12,6).Value = 583
$ws4.Cells.Item(14,6).Value = 3183
$ws4.Cells.Item(15,6).Value = 2759
$ws4.Cells.Item(16,6).Value = 532
$ws4.Cells.Item(19,6).Value = 355
$ws4.Cells.Item(20,6).Value = 338
$ws4.Cells.Item(21,6).Value = 62
$ws4.Cells.Item(22,6).Value = 291
$ws4.Cells.Item(27,6).Value = 635
$ws4.Cells.Item(31,6).Value = 465
$ws4.Cells.Item(38,6).Value = 292
$ws4.Cells.Item(39,6).Value = 1248
$ws4.Cells.Item(42,6).Value = 62
$ws4.Cells.Item(43,6).Value = 339
$ws4.Cells.Item(46,6).Value = 128
$ws4.Cells.Item(48,6).Value = 45
